$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's gold-price row (row 6): date in column A,
# price summary text in column B - matching the style of the existing
# data rows (2-4).
$ws.Range("A6").Value = "27-09-2025"
$ws.Range("B6").Value = "The price of gold in India today is ₹11,548 per gram for 24 karat gold, ₹10,585 per gram for 22 karat gold and ₹8,661 per gram for 18 karat gold (also called 999 gold)."
